# Apply the small "NS command output" corrections to the clan log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Bam) - donated units / ratio corrected
$ws.Range("G5").Value = 1276.0
$ws.Range("I5").Value = 0.64

# Row 7 (shag tand) - trophies / received units / ratio corrected
$ws.Range("E7").Value = 5066.0
$ws.Range("H7").Value = 1423.0
$ws.Range("I7").Value = 1.33

# Row 14 (GJS) - received units corrected
$ws.Range("H14").Value = 3674.0

# Rows 22 & 23 swap order: Bastos now outranks (j)de tik(j)
$ws.Range("B22").Value = "Bastos"
$ws.Range("C22").Value = "#8RP8QV8V"
$ws.Range("D22").Value = 172.0
$ws.Range("E22").Value = 4555.0
$ws.Range("F22").Value = "member"
$ws.Range("G22").Value = 37.0
$ws.Range("H22").Value = 0.0
$ws.Range("I22").Value = 37.0

$ws.Range("B23").Value = "(j)de tik(j)"
$ws.Range("C23").Value = "#GYVQ0Y8R"
$ws.Range("D23").Value = 178.0
$ws.Range("E23").Value = 4554.0
$ws.Range("F23").Value = "coLeader"
$ws.Range("G23").Value = 614.0
$ws.Range("H23").Value = 1624.0
$ws.Range("I23").Value = 0.38

# Row 28 (elandro) - trophies / received units / ratio corrected
$ws.Range("E28").Value = 4336.0
$ws.Range("H28").Value = 4001.0
$ws.Range("I28").Value = 0.79

# Row 34 (PeterClash) - received units corrected
$ws.Range("H34").Value = 5768.0

# Footer timestamp refined (08:32 -> 08:57)
$ws.PageSetup.LeftFooter = "Clanoverzicht"
$ws.PageSetup.RightFooter = "21/12/2017 08:57"
